# Generate Report for Handoff
# Adds two new tracked files (a .md file and a second .png dependency) to the
# localization-status report, renames the existing tracked file from a .md to
# a .png, and refreshes the handoff timestamps across all three sheets
# (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$e2eBase    = "https://github.com/OpenLocalizationTest/oltest/blob/f1fe55525e698ee636ba213cf394eb8b8ead2bf0/e2e/"
$zhBase     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09f5f96310a68c9dfe5b7dc6ca384d7d35884057/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deBase     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0d5dad7bc6340cff6bda51432a57e7486372a156/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$renamedFile   = "3d5f03c8-da09-422d-be3e-7bb1f7850cdc.png"
$mdFile        = "669e84ea-6871-438d-92a2-51c2ddafde38.md"
$pngFile       = "73099f27-6751-4323-8b59-585687972500.png"

$zhRenamedTarget = "6e263573f662e429c2d970904925a792eb159508.png"
$zhMdTarget      = "669e84ea-6871-438d-92a2-51c2ddafde38.8e8dccee38233a151eb3aff0ff1e4eb07d2f8d33.zh-cn.xlf"
$zhPngTarget     = "c662cf7d159b5c32c38efc4a0c59c862de48f896.png"

$deRenamedTarget = "6e263573f662e429c2d970904925a792eb159508.png"
$deMdTarget      = "669e84ea-6871-438d-92a2-51c2ddafde38.8e8dccee38233a151eb3aff0ff1e4eb07d2f8d33.de-de.xlf"
$dePngTarget     = "c662cf7d159b5c32c38efc4a0c59c862de48f896.png"

$readyStatus   = "Ready for handoff"
$overviewDate  = "2016-43-12 08:43:27"
$zeroDate      = "0001-01-01 00:00:00"
$zhHandoffDate = "2016-03-12 08:43:24"
$deHandoffDate = "2016-03-12 08:43:27"
$dependencyRef = "e2e\669e84ea-6871-438d-92a2-51c2ddafde38.md"

$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Set-CellHyperlink {
    param($ws, [string]$cellRef, [string]$address, [string]$display)
    $ws.Range($cellRef).Value = $display
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $display)
}

function Set-DateTextCell {
    param($ws, [string]$cellRef, [string]$text)
    $ws.Range($cellRef).NumberFormat = $dateFormat
    $ws.Range($cellRef).Value = $text
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Hyperlinks.Delete() removes every hyperlink on the sheet (not just the
# target range's), so clear them all up front and re-add every one below.
$wsOverview.Hyperlinks.Delete()

# Row 2: existing tracked file was renamed from .md to .png, handoff time bumped
Set-CellHyperlink $wsOverview "A2" ($e2eBase + $renamedFile) $renamedFile
$wsOverview.Range("D2").Value = $overviewDate

# Row 3 (new): the .md file that was added
Set-CellHyperlink $wsOverview "A3" ($e2eBase + $mdFile) $mdFile
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus
$wsOverview.Range("D3").Value = $overviewDate

# Row 4 (new): the second .png dependency file that was added
Set-CellHyperlink $wsOverview "A4" ($e2eBase + $pngFile) $pngFile
$wsOverview.Range("B4").Value = $readyStatus
$wsOverview.Range("C4").Value = $readyStatus
$wsOverview.Range("D4").Value = $overviewDate

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

# Row 2: renamed file, now flagged as a dependency (IsDependency + Dependency From)
Set-CellHyperlink $wsZh "A2" ($e2eBase + $renamedFile) $renamedFile
Set-CellHyperlink $wsZh "B2" ($e2eBase + $renamedFile) ".png"
$wsZh.Range("C2").Value = $readyStatus
Set-CellHyperlink $wsZh "D2" ($zhBase + $zhRenamedTarget) $zhRenamedTarget
Set-DateTextCell $wsZh "E2" $zhHandoffDate
$wsZh.Range("H2").Value = $zeroDate
$wsZh.Range("I2").Value = "IsDependency"
$wsZh.Range("J2").Value = $dependencyRef

# Row 3 (new): the .md file, included (not a dependency)
Set-CellHyperlink $wsZh "A3" ($e2eBase + $mdFile) $mdFile
Set-CellHyperlink $wsZh "B3" ($e2eBase + $mdFile) ".md"
$wsZh.Range("C3").Value = $readyStatus
Set-CellHyperlink $wsZh "D3" ($zhBase + $zhMdTarget) $zhMdTarget
Set-DateTextCell $wsZh "E3" $zhHandoffDate
$wsZh.Range("H3").Value = $zeroDate
$wsZh.Range("I3").Value = "Include"

# Row 4 (new): the second .png dependency file
Set-CellHyperlink $wsZh "A4" ($e2eBase + $pngFile) $pngFile
Set-CellHyperlink $wsZh "B4" ($e2eBase + $pngFile) ".png"
$wsZh.Range("C4").Value = $readyStatus
Set-CellHyperlink $wsZh "D4" ($zhBase + $zhPngTarget) $zhPngTarget
Set-DateTextCell $wsZh "E4" $zhHandoffDate
$wsZh.Range("H4").Value = $zeroDate
$wsZh.Range("I4").Value = "IsDependency"
$wsZh.Range("J4").Value = $dependencyRef

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

# Row 2: renamed file, now flagged as a dependency (IsDependency + Dependency From)
Set-CellHyperlink $wsDe "A2" ($e2eBase + $renamedFile) $renamedFile
Set-CellHyperlink $wsDe "B2" ($e2eBase + $renamedFile) ".png"
$wsDe.Range("C2").Value = $readyStatus
Set-CellHyperlink $wsDe "D2" ($deBase + $deRenamedTarget) $deRenamedTarget
Set-DateTextCell $wsDe "E2" $deHandoffDate
$wsDe.Range("H2").Value = $zeroDate
$wsDe.Range("I2").Value = "IsDependency"
$wsDe.Range("J2").Value = $dependencyRef

# Row 3 (new): the .md file, included (not a dependency)
Set-CellHyperlink $wsDe "A3" ($e2eBase + $mdFile) $mdFile
Set-CellHyperlink $wsDe "B3" ($e2eBase + $mdFile) ".md"
$wsDe.Range("C3").Value = $readyStatus
Set-CellHyperlink $wsDe "D3" ($deBase + $deMdTarget) $deMdTarget
Set-DateTextCell $wsDe "E3" $deHandoffDate
$wsDe.Range("H3").Value = $zeroDate
$wsDe.Range("I3").Value = "Include"

# Row 4 (new): the second .png dependency file
Set-CellHyperlink $wsDe "A4" ($e2eBase + $pngFile) $pngFile
Set-CellHyperlink $wsDe "B4" ($e2eBase + $pngFile) ".png"
$wsDe.Range("C4").Value = $readyStatus
Set-CellHyperlink $wsDe "D4" ($deBase + $dePngTarget) $dePngTarget
Set-DateTextCell $wsDe "E4" $deHandoffDate
$wsDe.Range("H4").Value = $zeroDate
$wsDe.Range("I4").Value = "IsDependency"
$wsDe.Range("J4").Value = $dependencyRef
